$wb = $excel.ActiveWorkbook

# --- ALC sheet updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 141.36842
$ws.Range("I6").Value = 83.76922999999999
$ws.Range("J6").Value = 266.16666
$ws.Range("K6").Value = 251.30769
$ws.Range("L6").Value = 798.4999799999999
$ws.Range("M6").Value = -139.30769
$ws.Range("N6").Value = -1022.49998
$ws.Range("H11").Value = 142857150
$ws.Range("I11").Value = 142857150
$ws.Range("K11").Value = 142857150
$ws.Range("M11").Value = -142857010
$ws.Range("H86").Value = 7650757.5
$ws.Range("I86").Value = 100000
$ws.Range("J86").Value = 15201515
$ws.Range("K86").Value = 100000
$ws.Range("L86").Value = 15201515
$ws.Range("M86").Value = -98877
$ws.Range("N86").Value = -15203761
$ws.Range("H89").Value = 7650757.5
$ws.Range("I89").Value = 100000
$ws.Range("J89").Value = 15201515
$ws.Range("K89").Value = 500000
$ws.Range("L89").Value = 76007575
$ws.Range("M89").Value = -494384
$ws.Range("N89").Value = -76018807
$ws.Range("H106").Value = 3054
$ws.Range("I106").Value = 2246.6667
$ws.Range("K106").Value = 2246.6667
$ws.Range("M106").Value = -1615.6667

# --- ARM sheet updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 35745.93
$ws.Range("I45").Value = 67728.92999999999
$ws.Range("K45").Value = 67728.92999999999
$ws.Range("M45").Value = -67351.92999999999
$ws.Range("H122").Value = 1073.1428
$ws.Range("I122").Value = 1122.4
$ws.Range("J122").Value = 950
$ws.Range("K122").Value = 3367.2
$ws.Range("L122").Value = 2850
$ws.Range("M122").Value = -917.2000000000003
$ws.Range("N122").Value = -7750

# --- CRP sheet updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2520.2778
$ws.Range("I16").Value = 1121.5
$ws.Range("J16").Value = 3639.3
$ws.Range("K16").Value = 1121.5
$ws.Range("L16").Value = 3639.3
$ws.Range("M16").Value = -834.5
$ws.Range("N16").Value = -4213.3
$ws.Range("H22").Value = 15914.714
$ws.Range("I22").Value = 275.75
$ws.Range("J22").Value = 36766.668
$ws.Range("K22").Value = 275.75
$ws.Range("L22").Value = 36766.668
$ws.Range("M22").Value = 74.25
$ws.Range("N22").Value = -37466.668
$ws.Range("H113").Value = 2520.2778
$ws.Range("I113").Value = 1121.5
$ws.Range("J113").Value = 3639.3
$ws.Range("K113").Value = 1121.5
$ws.Range("L113").Value = 3639.3
$ws.Range("M113").Value = 1048.5
$ws.Range("N113").Value = -7979.3
$ws.Range("H122").Value = 35715176
$ws.Range("I122").Value = 35715176
$ws.Range("K122").Value = 107145528
$ws.Range("M122").Value = -107143078
$ws.Range("H141").Value = 33343.4
$ws.Range("J141").Value = 33343.4
$ws.Range("L141").Value = 33343.4
$ws.Range("N141").Value = -43703.4

# --- CUL sheet updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 949.75
$ws.Range("I35").Value = 999
$ws.Range("J35").Value = 933.3333
$ws.Range("K35").Value = 2997
$ws.Range("L35").Value = 2799.9999
$ws.Range("M35").Value = -2709
$ws.Range("N35").Value = -3375.9999
$ws.Range("H98").Value = 3002
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 3002
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 9006
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -12002
$ws.Range("H113").Value = 30117386
$ws.Range("I113").Value = 10417054
$ws.Range("J113").Value = 44444900
$ws.Range("K113").Value = 31251162
$ws.Range("L113").Value = 133334700
$ws.Range("M113").Value = -31248992
$ws.Range("N113").Value = -133339040
$ws.Range("H131").Value = 915.14
$ws.Range("I131").Value = 502
$ws.Range("J131").Value = 936.8842
$ws.Range("K131").Value = 1506
$ws.Range("L131").Value = 2810.6526
$ws.Range("M131").Value = 3534
$ws.Range("N131").Value = -12890.6526

# --- GSM sheet updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4867.1577
$ws.Range("I70").Value = 4534.8335
$ws.Range("J70").Value = 5436.857
$ws.Range("K70").Value = 4534.8335
$ws.Range("L70").Value = 5436.857
$ws.Range("M70").Value = -4264.8335
$ws.Range("N70").Value = -5976.857
$ws.Range("H73").Value = 4867.1577
$ws.Range("I73").Value = 4534.8335
$ws.Range("J73").Value = 5436.857
$ws.Range("K73").Value = 4534.8335
$ws.Range("L73").Value = 5436.857
$ws.Range("M73").Value = -3598.8335
$ws.Range("N73").Value = -7308.857
$ws.Range("H80").Value = 20005686
$ws.Range("I80").Value = 8544.333000000001
$ws.Range("K80").Value = 8544.333000000001
$ws.Range("M80").Value = -7546.333000000001
$ws.Range("H83").Value = 20005686
$ws.Range("I83").Value = 8544.333000000001
$ws.Range("K83").Value = 42721.665
$ws.Range("M83").Value = -37729.665
$ws.Range("H113").Value = 1400
$ws.Range("I113").Value = 1283.3334
$ws.Range("J113").Value = 1925
$ws.Range("K113").Value = 1283.3334
$ws.Range("L113").Value = 1925
$ws.Range("M113").Value = 886.6666
$ws.Range("N113").Value = -6265

# --- LTW sheet updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 2502.5
$ws.Range("I9").Value = 170
$ws.Range("J9").Value = 9500
$ws.Range("K9").Value = 170
$ws.Range("L9").Value = 9500
$ws.Range("M9").Value = 54
$ws.Range("N9").Value = -9948
$ws.Range("H40").Value = 8929993
$ws.Range("I40").Value = 1446.9524
$ws.Range("J40").Value = 35715628
$ws.Range("K40").Value = 1446.9524
$ws.Range("L40").Value = 35715628
$ws.Range("M40").Value = -1310.9524
$ws.Range("N40").Value = -35715900
$ws.Range("H63").Value = 4443.6665
$ws.Range("J63").Value = 4443.6665
$ws.Range("L63").Value = 4443.6665
$ws.Range("N63").Value = -5941.6665
$ws.Range("H66").Value = 4443.6665
$ws.Range("J66").Value = 4443.6665
$ws.Range("L66").Value = 13330.9995
$ws.Range("N66").Value = -20818.9995

# --- WVR sheet updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 9181.5
$ws.Range("I74").Value = 13000
$ws.Range("J74").Value = 7908.6665
$ws.Range("K74").Value = 13000
$ws.Range("L74").Value = 7908.6665
$ws.Range("M74").Value = -12064
$ws.Range("N74").Value = -9780.666499999999
$ws.Range("H77").Value = 9181.5
$ws.Range("I77").Value = 13000
$ws.Range("J77").Value = 7908.6665
$ws.Range("K77").Value = 39000
$ws.Range("L77").Value = 23725.9995
$ws.Range("M77").Value = -34320
$ws.Range("N77").Value = -33085.99950000001
$ws.Range("H119").Value = 33299.5
$ws.Range("J119").Value = 33299.5
$ws.Range("L119").Value = 33299.5
$ws.Range("N119").Value = -42975.5
$ws.Range("H122").Value = 30412.79
$ws.Range("J122").Value = 4719.3335
$ws.Range("L122").Value = 14158.0005
$ws.Range("N122").Value = -19058.0005
